# "actualizacion 19 de mayo"
#
# 1) LAURA ORTEGON sheet (5th tab): update the IMPUGNACIONES case tracking
#    (row 4) with new dated notes, and overwrite the old PUTUMAYO
#    revocatoria case (row 5) with a duplicate of the row below it
#    (CNE-E-DG-2025-001116), matching the author's edit.
# 2) Move the "active sheet" focus from JHON TRUJILLO to LAURA ORTEGON,
#    updating each sheet's remembered selection.

$wb = $excel.ActiveWorkbook

# ---- JHON TRUJILLO (tab 3): was the active tab; just leave behind a new
# remembered selection, no longer the active tab. ----
$wsTrujillo = $wb.Worksheets.Item(3)
$wsTrujillo.Range("P8").Select()

# ---- LAURA ORTEGON (tab 5) ----
$wsLaura = $wb.Worksheets.Item(5)

# Row 4: CNE-E-DG-2023-014262 - CNE-E-DG-2023-014265 (IMPUGNACIONES / GUAJIRA)
# TRAZABILIDAD (J4): new entry prepended
$oldJ4 = $wsLaura.Range("J4").Text
$wsLaura.Range("J4").Value = "01/04/2025 REOSLUCION 01459 QUE CONCEDE IMPUGNACION`n" + $oldJ4

# FECHA DE ULTIMA ACTUACIÓN (L4): 16/01/2025 -> 01/04/2025
$wsLaura.Range("L4").Value = 45748

# SEGUIMIENTO (M4): new entry prepended
$oldM4 = $wsLaura.Range("M4").Text
$wsLaura.Range("M4").Value = "06/05/2025 RECIBIDA CONSTANCIA DE NOTIFICACION DE LA REOSLUCION 01459 QUE CONCEDE IMPUGNACION`n" + $oldM4

# FECHA DE SEGUIMIENTO (N4): 08/04/2025 -> 06/05/2025
$wsLaura.Range("N4").Value = 45783

# Row grew by one more line of wrapped text -> taller row
$wsLaura.Rows.Item(4).RowHeight = 369.75

# Row 5 (CNE-E-DG-2025-000141, PUTUMAYO revocatoria) gets overwritten with a
# copy of row 6 (CNE-E-DG-2025-001116) -- formats first, then values, mirroring
# a copy/paste of the row above onto it.
$srcRow = $wsLaura.Range("A6:Q6")
$dstRow = $wsLaura.Range("A5:Q5")
$srcRow.Copy()
$dstRow.PasteSpecial(-4122) # xlPasteFormats
$srcRow.Copy()
$dstRow.PasteSpecial(-4163) # xlPasteValues
$excel.CutCopyMode = $false

# LAURA ORTEGON becomes the active tab, with A2:Q5 highlighted.
$wsLaura.Activate()
$wsLaura.Range("A2:Q5").Select()
